# Generate Report for Handoff
# Adds a new handed-off file ("ffe9c94c-...") as row 3 in the Overview,
# zh-cn and de-de tables, mirroring the existing row 2 pattern.

$wb = $excel.ActiveWorkbook

$mdName    = "ffe9c94c-ad82-4d9b-b864-a3160e651f9aooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$e2eMdName = "e2e\ffe9c94c-ad82-4d9b-b864-a3160e651f9aooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$ready     = "Ready for handoff"
$dateHo    = "2016-09-01 08:34:51"
$zhXlf     = "ffe9c94c-ad82-4d9b-b864-a3160e651f9aoooooooooooooooooooooooooooooooooooooooo.65e0c7b3b11b3251ae4f94b7295f756ee460895f.zh-cn.xlf"
$dateZh    = "2016-09-01 08:34:47"
$deXlf     = "ffe9c94c-ad82-4d9b-b864-a3160e651f9aoooooooooooooooooooooooooooooooooooooooo.65e0c7b3b11b3251ae4f94b7295f756ee460895f.de-de.xlf"
$zeroDate  = "0001-01-01 00:00:00"
$hyperUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ced822f7cc9f242d7225759ee8941f7d1c1eb4b9/e2e/ffe9c94c-ad82-4d9b-b864-a3160e651f9aooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
# Leading apostrophe forces these to be stored as literal text instead of
# being auto-coerced to Boolean TRUE/FALSE by the COM value setter.
$trueTxt   = "'True"
$falseTxt  = "'False"

# ---------------------------------------------------------------------------
# Overview sheet: append row 3 to the "Overview" table
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add()

$wsOverview.Range("A3").Value = $mdName
$wsOverview.Range("B3").Value = $e2eMdName
$wsOverview.Range("B3").Style = "HyperLink"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperUrl, "", "", $e2eMdName)
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $ready
$wsOverview.Range("F3").Value = $ready
$wsOverview.Range("G3").Value = $dateHo
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Columns.Item(5).ColumnWidth = 16.35
$wsOverview.Columns.Item(6).ColumnWidth = 16.35

# ---------------------------------------------------------------------------
# zh-cn sheet: append row 3 to the "zh-cn" table
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add()

$wsZhCn.Range("A3").Value = $mdName
$wsZhCn.Range("A3").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperUrl, "", "", $mdName)
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $ready
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = $falseTxt
$wsZhCn.Range("G3").Value = $zhXlf
$wsZhCn.Range("H3").Value = $dateZh
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = $zeroDate
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = $trueTxt
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = $falseTxt
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = 16.35

# ---------------------------------------------------------------------------
# de-de sheet: append row 3 to the "de-de" table
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add()

$wsDeDe.Range("A3").Value = $mdName
$wsDeDe.Range("A3").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperUrl, "", "", $mdName)
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $ready
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = $falseTxt
$wsDeDe.Range("G3").Value = $deXlf
$wsDeDe.Range("H3").Value = $dateHo
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = $zeroDate
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = $trueTxt
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = $falseTxt
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 16.35

Write-Host "Done"
